$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Renaming an InlineShape's docPr/name needs to go through the current
# Selection (Shape.Select() + $word.Selection.InlineShapes) - calling
# .Name directly on a handle obtained from a Footer range can otherwise
# report a stale-handle error, so route every rename through Select().
function Rename-PictureShape($shape, $newName) {
    $shape.Select()
    $word.Selection.InlineShapes.Item(1).Name = $newName
}

# Section 1 only has "default" (Item 1) and "first page" (Item 2)
# headers/footers - this document has no even-page header/footer.

# header1.xml (BTec logo, docPr id="1") -> Headers.Item(2) ("first page")
Rename-PictureShape $sec.Headers.Item(2).Range.InlineShapes.Item(1) "image1.jpg"

# header2.xml (BTec logo, docPr id="3") -> Headers.Item(1) ("default")
Rename-PictureShape $sec.Headers.Item(1).Range.InlineShapes.Item(1) "image1.jpg"

# footer1.xml (Pearson logo, docPr id="2") -> Footers.Item(2) ("first page")
Rename-PictureShape $sec.Footers.Item(2).Range.InlineShapes.Item(1) "image2.png"

# footer2.xml (Pearson logo, docPr id="4") -> Footers.Item(1) ("default")
Rename-PictureShape $sec.Footers.Item(1).Range.InlineShapes.Item(1) "image2.png"
